# Generate Report for Handoff
# Inserts a new row (for the file
# "active-directory-aadconnectsync-understanding-declarative-provisioning-expressions.md")
# above the existing "test-content-1.md" row on every sheet (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$newFileBase = "active-directory-aadconnectsync-understanding-declarative-provisioning-expressions"
$newFileName = "$newFileBase.md"
$newCommit   = "be82fa1f69efd7bef6381c68ef71173a98cd4292"
$newXlfZh    = "$newFileBase.$newCommit.zh-cn.xlf"
$newXlfDe    = "$newFileBase.$newCommit.de-de.xlf"

$mdBaseUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/298395042849475f4d5e9297757af5e7954e65d5/e2e"
$cfgUrl      = "https://github.com/OpenLocalizationTest/oltest/blob/298395042849475f4d5e9297757af5e7954e65d5/.localization-config"
$xlfZhBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f21e7fa7ecd381ccd9e49904b515d65b392b6e0f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho"
$xlfDeBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/26df1a5acbbe3ac5fd71a9087d4f30eba0766425/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho"

# ---------------------------------------------------------------------------
# Sheet "Overview": columns A (File Name) / B (zh-cn status) / C (de-de status)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Rows(2).Insert()

$ws.Range("A2").Value = $newFileName
$ws.Range("A2").Style = "HyperLink"
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "$mdBaseUrl/$newFileName", "", "", $newFileName) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "$mdBaseUrl/test-content-1.md", "", "", "test-content-1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn": full status table
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Rows(2).Insert()

$ws.Range("A2").Value = $newFileName
$ws.Range("A2").Style = "HyperLink"
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = $newXlfZh
$ws.Range("C2").Style = "HyperLink"
$ws.Range("D2").Value = "2016-02-17 05:26:14"
$ws.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G2").Value = "0001-01-01 00:00:00"
$ws.Range("H2").Value = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "$mdBaseUrl/$newFileName", "", "", $newFileName) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "$xlfZhBase/$newXlfZh", "", "", $newXlfZh) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "$mdBaseUrl/test-content-1.md", "", "", "test-content-1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "$xlfZhBase/test-content-1.fc92bd3ea58825d16dce72b4174897baa7b370e7.zh-cn.xlf", "", "", "test-content-1.fc92bd3ea58825d16dce72b4174897baa7b370e7.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de": full status table
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows(2).Insert()

$ws.Range("A2").Value = $newFileName
$ws.Range("A2").Style = "HyperLink"
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = $newXlfDe
$ws.Range("C2").Style = "HyperLink"
$ws.Range("D2").Value = "2016-02-17 05:26:24"
$ws.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G2").Value = "0001-01-01 00:00:00"
$ws.Range("H2").Value = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "$mdBaseUrl/$newFileName", "", "", $newFileName) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "$xlfDeBase/$newXlfDe", "", "", $newXlfDe) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "$mdBaseUrl/test-content-1.md", "", "", "test-content-1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "$xlfDeBase/test-content-1.fc92bd3ea58825d16dce72b4174897baa7b370e7.de-de.xlf", "", "", "test-content-1.fc92bd3ea58825d16dce72b4174897baa7b370e7.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, "", "", ".localization-config") | Out-Null

Write-Output "Report for handoff generated."
